# Refresh the cryptos price list (GitHub Actions scheduled update).
# Only B/C/D/E of the data rows change; the row-index column A is untouched.
# Some new "Price" strings are plain decimal numbers (e.g. "245.04") which
# Excel would otherwise auto-coerce to a Number on assignment; a leading
# apostrophe (quote-prefix) keeps those as Text, matching the original
# inline-string cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "35.503.56"
$ws.Range("E2").Value = "  -2.75%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.977.11"
$ws.Range("E3").Value = "  -3.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "245.04"
$ws.Range("E5").Value = "  +1.14%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -4.65%  "

# Row 7 - Solana
Set-TextValue "D7" "56.99"
$ws.Range("E7").Value = "  +4.38%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - OKB
Set-TextValue "D9" "58.40"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.56%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0731"
$ws.Range("E11").Value = "  -2.48%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -2.80%  "

# Row 13 - Polygon
Set-TextValue "D13" "0.939"
$ws.Range("E13").Value = "  +3.35%  "

# Row 14 - Chainlink
Set-TextValue "D14" "14.40"
$ws.Range("E14").Value = "  -2.29%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.263.70"

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -2.75%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.007.74"
$ws.Range("E17").Value = "  -2.18%  "

# Row 18 - Avalanche
Set-TextValue "D18" "17.83"
$ws.Range("E18").Value = "  +6.33%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "35.498.82"
$ws.Range("E19").Value = "  -2.56%  "

# Row 20 - Litecoin
Set-TextValue "D20" "71.27"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -1.97%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "232.51"
$ws.Range("E22").Value = "  -2.43%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.82%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.01%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.50"
$ws.Range("E25").Value = "  +17.99%  "

# Row 26 - Toncoin
Set-TextValue "D26" "2.29"
$ws.Range("E26").Value = "  -2.71%  "

# Row 27 - Monero
Set-TextValue "D27" "164.20"
$ws.Range("E27").Value = "  -0.32%  "

# Row 28 - Cosmos
Set-TextValue "D28" "9.06"
$ws.Range("E28").Value = "  -3.08%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.07"
$ws.Range("E29").Value = "  -4.90%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  -2.84%  "

# Row 31 - Filecoin
Set-TextValue "D31" "4.84"
$ws.Range("E31").Value = "  -4.59%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -9.46%  "

# Row 33 - Kaspa
Set-TextValue "D33" "0.0943"
$ws.Range("E33").Value = "  +15.26%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0589"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - now LidoDAOToken (was InternetComputer(DFINITY))
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D35" "2.39"
$ws.Range("E35").Value = "  +8.27%  "

# Row 36 - now InternetComputer(DFINITY) (was LidoDAOToken)
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D36" "4.29"
$ws.Range("E36").Value = "  -3.54%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.10%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -4.42%  "

# Row 39 - THORChain
Set-TextValue "D39" "5.16"
$ws.Range("E39").Value = "  +6.77%  "

# Row 40 - TrustWalletToken
$ws.Range("E40").Value = "  -2.53%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +1.98%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  -2.64%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  -1.85%  "

# Row 44 - Aave
Set-TextValue "D44" "91.20"
$ws.Range("E44").Value = "  -2.74%  "

# Row 45 - InjectiveProtocol
Set-TextValue "D45" "15.89"
$ws.Range("E45").Value = "  -0.08%  "

# Row 46 - now FraxShare (was Cronos)
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "7.50"
$ws.Range("E46").Value = "  -0.84%  "

# Row 47 - now Cronos (was FraxShare)
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.0880"
$ws.Range("E47").Value = "  -6.61%  "

# Row 48 - Maker
Set-TextValue "D48" "1.366.49"
$ws.Range("E48").Value = "  -3.26%  "

# Row 49 - MXToken
Set-TextValue "D49" "2.89"
$ws.Range("E49").Value = "  +1.19%  "

# Row 50 - MultiversX
Set-TextValue "D50" "46.78"
$ws.Range("E50").Value = "  +3.81%  "

# Row 51 - FTXToken
Set-TextValue "D51" "3.65"
$ws.Range("E51").Value = "  +12.47%  "
